$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows at 13:14 -- this shifts old rows 13-23 down to 15-25,
#    carrying their column-A labels & formatting along automatically.
$ws.Range("A13:A14").EntireRow.Insert()

# 2. The inherited row height for (old row 18 -> new) row 20 is wrong (60); restore default.
$ws.Rows.Item(20).RowHeight = 15

# 3. Remove the stray empty styled A-cells that Insert() left behind on the two new rows.
$ws.Range("A13:A14").ClearContents()
$ws.Range("A13:A14").ClearFormats()

# 4. Populate the two new rows (B/C) with the professor names, copying the wrap styles
#    (style index 2 for column B, style index 3 for column C) from an existing template row.
$ws.Cells.Item(3,2).Copy($ws.Cells.Item(13,2))
$ws.Cells.Item(3,3).Copy($ws.Cells.Item(13,3))
$ws.Cells.Item(3,2).Copy($ws.Cells.Item(14,2))
$ws.Cells.Item(3,3).Copy($ws.Cells.Item(14,3))

$ws.Cells.Item(13,2).Value = '6495737 - Durval Rodrigues Junior'
$ws.Cells.Item(13,3).Value = '6495737 - Durval Rodrigues Junior'
$ws.Cells.Item(14,2).Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Cells.Item(14,3).Value = '1643715 - Paulo Atsushi Suzuki'

# 5. Correct the body text in the rows whose B/C content was wrong (duplicated from other
#    rows in the 'before' file) -- styles/heights on these rows are already correct after the
#    row insert & shift, only the text needs updating.
# Row 10: Objetivos / Objectives body text
$ws.Cells.Item(10,2).Value = 'Apresentação introdutória das técnicas de análise microestrutural de materiais. Apresentação das técnicas e equipamentos necessários para a análise microestrutural. Seleção adequada das técnicas experimentais. Apresentação das técnicas adequadas de preparação de amostras. Verificação dos custos envolvidos nas técnicas de caracterização microestrutural.'
$ws.Cells.Item(10,3).Value = 'Apresentação introdutória das técnicas de análise microestrutural de materiais. Apresentação das técnicas e equipamentos necessários para a análise microestrutural. Seleção adequada das técnicas experimentais. Apresentação das técnicas adequadas de preparação de amostras. Verificação dos custos envolvidos nas técnicas de caracterização microestrutural.'

# Row 15: Programa resumido / Short syllabus body text
$ws.Cells.Item(15,2).Value = 'A Microestrutura dos Materiais. Difratometria de raios X. Análise Microestrutural utilizando Luz Síncrotron. Microscopia Óptica. Microscopia Eletrônica. Microscopia de Tunelamento e de Força Atômica. Análise Química de Microrregiões. Análises Térmicas. Fluorescência de raios X. Técnicas Indiretas de Análise de Microestrutura. Seleção de Técnicas Experimentais.'
$ws.Cells.Item(15,3).Value = 'A Microestrutura dos Materiais. Difratometria de raios X. Análise Microestrutural utilizando Luz Síncrotron. Microscopia Óptica. Microscopia Eletrônica. Microscopia de Tunelamento e de Força Atômica. Análise Química de Microrregiões. Análises Térmicas. Fluorescência de raios X. Técnicas Indiretas de Análise de Microestrutura. Seleção de Técnicas Experimentais.'

# Row 17: Programa / Syllabus body text
$ws.Cells.Item(17,2).Value = '1. A Microestrutura dos Materiais. 2. Difratometria de raios X. 3. Análise Microestrutural utilizando Luz Síncrotron. 4. Microscopia Óptica. 5. Microscopia Eletrônica. 6. Microscopia de Tunelamento e de Força Atômica. 7. Análise Química de Microrregiões. 8. Análises Térmicas. 9. Fluorescência de raios X. 10. Técnicas Indiretas de Análise de Microestrutura. 11. Seleção de Técnicas Experimentais.'
$ws.Cells.Item(17,3).Value = '1. A Microestrutura dos Materiais. 2. Difratometria de raios X. 3. Análise Microestrutural utilizando Luz Síncrotron. 4. Microscopia Óptica. 5. Microscopia Eletrônica. 6. Microscopia de Tunelamento e de Força Atômica. 7. Análise Química de Microrregiões. 8. Análises Térmicas. 9. Fluorescência de raios X. 10. Técnicas Indiretas de Análise de Microestrutura. 11. Seleção de Técnicas Experimentais.'

# Row 20: Método body text
$ws.Cells.Item(20,2).Value = 'Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre).'
$ws.Cells.Item(20,3).Value = 'Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre).'

# Row 21: Critério body text
$ws.Cells.Item(21,2).Value = 'A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2.'
$ws.Cells.Item(21,3).Value = 'A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2.'

# Row 22: Norma de recuperação body text
$ws.Cells.Item(22,2).Value = 'Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação).'
$ws.Cells.Item(22,3).Value = 'Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação).'

# Row 23: Bibliografia body text
$ws.Cells.Item(23,2).Value = '1. Van Vlack, L.H. Princípios de Ciência e Tecnologia dos Materiais, 4a.ed., Ed. Campus, Rio de Janeiro, 1984. 2. Shackelford, J.F. Introduction to Materials Science for Engineers. 4th Edition. Prentice Hall Inc., 1996. 3. Padilha, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985. 4. Guy, A.G. Ciência dos Materiais. Livros Técnicos e Científicos Editora, 1982. 5. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 6. Nondestructive Characterization of Materials. Series. Plenum Press, New York. 7. Yacobi, B.G. Holt, D.B. Kazmerski, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994. 8. Lowell, S.; Shields, J. E.; Thomas, M. A.; Thommes, M. Characterization of Porous Solids and Powders: Surface Area, Pore Size and Density, Springer, 2010. 9. Murphy, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001. 10. Wu, Q.; Merchant, F.; Castleman, K. Microscope Image Processing, Academic Press, 2008. 11. Cullity, B. D.; Stock, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001. 12. Goldstein, J.; et al., Scanning Electron Microscopy and X-ray Microanalysis, Springer, 2003. 13. Hatakeyama, T.; Zhenhai, L. Handbook of Thermal Analysis, NY: Wiley, 1999. 14. Haines, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002. 15. Schramm, G. Reologia e Reometria. Editora Artliber, 2006.16. Azevedo, A. D.; Mothe, C. G. Análise Térmica de Materiais. São Paulo: ARTLIBER, 2009.17. Brown, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.18. Muller, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.19. Speyer, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.'
$ws.Cells.Item(23,3).Value = '1. Van Vlack, L.H. Princípios de Ciência e Tecnologia dos Materiais, 4a.ed., Ed. Campus, Rio de Janeiro, 1984. 2. Shackelford, J.F. Introduction to Materials Science for Engineers. 4th Edition. Prentice Hall Inc., 1996. 3. Padilha, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985. 4. Guy, A.G. Ciência dos Materiais. Livros Técnicos e Científicos Editora, 1982. 5. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 6. Nondestructive Characterization of Materials. Series. Plenum Press, New York. 7. Yacobi, B.G. Holt, D.B. Kazmerski, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994. 8. Lowell, S.; Shields, J. E.; Thomas, M. A.; Thommes, M. Characterization of Porous Solids and Powders: Surface Area, Pore Size and Density, Springer, 2010. 9. Murphy, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001. 10. Wu, Q.; Merchant, F.; Castleman, K. Microscope Image Processing, Academic Press, 2008. 11. Cullity, B. D.; Stock, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001. 12. Goldstein, J.; et al., Scanning Electron Microscopy and X-ray Microanalysis, Springer, 2003. 13. Hatakeyama, T.; Zhenhai, L. Handbook of Thermal Analysis, NY: Wiley, 1999. 14. Haines, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002. 15. Schramm, G. Reologia e Reometria. Editora Artliber, 2006.16. Azevedo, A. D.; Mothe, C. G. Análise Térmica de Materiais. São Paulo: ARTLIBER, 2009.17. Brown, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.18. Muller, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.19. Speyer, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.'

